# Applies the "Penalty Reward System" (unfinished) edits:
#  - Forecast Comparison sheet: shift each row's Week_Start_Date (col B) to the
#    next week's date, and update MyForecast (col D) values.
#  - Summary sheet: update several derived metrics to match the new data.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: rows 2-17 -----------------------------
# row, new Week_Start_Date (col B), new MyForecast (col D)
$rows = @(
    @(2,  "2025-01-12", 0),
    @(3,  "2025-01-19", 0),
    @(4,  "2025-01-26", 0),
    @(5,  "2025-02-02", 0),
    @(6,  "2025-02-09", 0),
    @(7,  "2025-02-16", 0),
    @(8,  "2025-02-23", 1),
    @(9,  "2025-03-02", 0),
    @(10, "2025-03-09", 0),
    @(11, "2025-03-16", 1),
    @(12, "2025-03-23", 1),
    @(13, "2025-03-30", 1),
    @(14, "2025-04-06", 1),
    @(15, "2025-04-13", 1),
    @(16, "2025-04-20", 1),
    @(17, "2025-04-27", 1)
)

foreach ($row in $rows) {
    $r = $row[0]
    $newDate = $row[1]
    $newForecast = $row[2]

    # Leading apostrophe forces Excel to keep the value as literal text
    # instead of re-interpreting the date-looking string as a date serial.
    $wsForecast.Cells.Item($r, 2).Value = "'" + $newDate
    $wsForecast.Cells.Item($r, 4).Value = $newForecast
}

# --- Summary sheet -----------------------------------------------------
$wsSummary.Range("B2").Value = "2022-12-25 to 2025-01-05"

$wsSummary.Range("B9").Value = "'9"
$wsSummary.Range("B10").Value = "'4"
$wsSummary.Range("B11").Value = "'2"
$wsSummary.Range("B12").Value = "'1"
$wsSummary.Range("B13").Value = "'2025-04-20"
$wsSummary.Range("B14").Value = "'0"
$wsSummary.Range("B15").Value = "'2025-01-12"
